$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.001257564406841993
$ws.Range("E2").Value = 0.5737204579636455
$ws.Range("G2").Value = 0.03989016404375434
$ws.Range("H2").Value = 0.2555096605792642
$ws.Range("I2").Value = 0.09472202695906162
$ws.Range("J2").Value = 0.1415524398908019
$ws.Range("K2").Value = 0.01031477889046073
$ws.Range("D3").Value = 0.07872134773060679
$ws.Range("E3").Value = 0.5977131072431803
$ws.Range("G3").Value = 0.0345823522657156
$ws.Range("H3").Value = 0.3188638077117503
$ws.Range("I3").Value = 0.05538225639611483
$ws.Range("J3").Value = 0.1505562500096858
$ws.Range("K3").Value = 0.00933563569560647
$ws.Range("C4").Value = 750
$ws.Range("D4").Value = 0.08191800210624933
$ws.Range("E4").Value = 0.6067328019998968
$ws.Range("F4").Value = 750
$ws.Range("G4").Value = 0.03412511944770813
$ws.Range("H4").Value = 0.3267448549158871
$ws.Range("I4").Value = 0.05122113320976496
$ws.Range("J4").Value = 0.1559512559324503
$ws.Range("K4").Value = 0.009595633018761873
$ws.Range("D5").Value = 0.001803758088499308
$ws.Range("E5").Value = 0.5769381779246032
$ws.Range("G5").Value = 0.03891843976452947
$ws.Range("H5").Value = 0.2565065808594227
$ws.Range("I5").Value = 0.09635335579514503
$ws.Range("J5").Value = 0.1425050417892635
$ws.Range("K5").Value = 0.01048250822350383
$ws.Range("C6").Value = 2399
$ws.Range("D6").Value = 0.1788086406886578
$ws.Range("E6").Value = 3.488423581700772
$ws.Range("F6").Value = 2399
$ws.Range("G6").Value = 0.1105686975643039
$ws.Range("H6").Value = 0.9602491110563278
$ws.Range("I6").Value = 1.979463193099946
$ws.Range("J6").Value = 0.3030690466985106
$ws.Range("K6").Value = 0.03880238952115178
$ws.Range("E7").Value = 7.253541064914316
$ws.Range("D8").Value = 0.000702104065567255
$ws.Range("E8").Value = 0.3880765549838543
$ws.Range("G8").Value = 0.0256541813723743
$ws.Range("H8").Value = 0.1828547627665102
$ws.Range("I8").Value = 0.05925528146326542
$ws.Range("J8").Value = 0.09035487053915858
$ws.Range("K8").Value = 0.007394487969577312
$ws.Range("D9").Value = 0.05967893823981285
$ws.Range("E9").Value = 0.4905664888210595
$ws.Range("G9").Value = 0.02645971858873963
$ws.Range("H9").Value = 0.2611825908534229
$ws.Range("I9").Value = 0.04930277541279793
$ws.Range("J9").Value = 0.1213390715420246
$ws.Range("K9").Value = 0.007955154869705439
$ws.Range("C10").Value = 702
$ws.Range("D10").Value = 0.05921263713389635
$ws.Range("E10").Value = 0.4920201976783574
$ws.Range("F10").Value = 702
$ws.Range("G10").Value = 0.02764697605744004
$ws.Range("H10").Value = 0.2645076452754438
$ws.Range("I10").Value = 0.04013616824522614
$ws.Range("J10").Value = 0.1270410129800439
$ws.Range("K10").Value = 0.008069343399256468
$ws.Range("D11").Value = 0.001043336000293493
$ws.Range("E11").Value = 0.388568548951298
$ws.Range("G11").Value = 0.02577080950140953
$ws.Range("H11").Value = 0.1810599220916629
$ws.Range("I11").Value = 0.06152039766311646
$ws.Range("J11").Value = 0.09039383241906762
$ws.Range("K11").Value = 0.0074429283849895
$ws.Range("D12").Value = 0.148574466817081
$ws.Range("E12").Value = 3.894684855360538
$ws.Range("G12").Value = 0.1143656615167856
$ws.Range("H12").Value = 0.9383131782524288
$ws.Range("I12").Value = 2.40365110617131
$ws.Range("J12").Value = 0.3003866509534419
$ws.Range("K12").Value = 0.03914503287523985
$ws.Range("E13").Value = 6.552316022105515
$ws.Range("D14").Value = 0.0009798342362046242
$ws.Range("E14").Value = 0.4808985064737499
$ws.Range("G14").Value = 0.03277217270806432
$ws.Range("H14").Value = 0.2191822116728872
$ws.Range("I14").Value = 0.07698865421116352
$ws.Range("J14").Value = 0.1159536552149802
$ws.Range("K14").Value = 0.008854633430019021
$ws.Range("D15").Value = 0.06920014298520982
$ws.Range("E15").Value = 0.5441397980321199
$ws.Range("G15").Value = 0.03052103542722762
$ws.Range("H15").Value = 0.2900231992825866
$ws.Range("I15").Value = 0.05234251590445638
$ws.Range("J15").Value = 0.1359476607758552
$ws.Range("K15").Value = 0.008645395282655954
$ws.Range("C16").Value = 726
$ws.Range("D16").Value = 0.07056531962007284
$ws.Range("E16").Value = 0.5493764998391271
$ws.Range("F16").Value = 726
$ws.Range("G16").Value = 0.03088604775257409
$ws.Range("H16").Value = 0.2956262500956655
$ws.Range("I16").Value = 0.04567865072749555
$ws.Range("J16").Value = 0.1414961344562471
$ws.Range("K16").Value = 0.008832488209009171
$ws.Range("D17").Value = 0.0014235470443964
$ws.Range("E17").Value = 0.4827533634379506
$ws.Range("G17").Value = 0.0323446246329695
$ws.Range("H17").Value = 0.2187832514755428
$ws.Range("I17").Value = 0.07893687672913074
$ws.Range("J17").Value = 0.1164494371041656
$ws.Range("K17").Value = 0.008962718304246664
$ws.Range("C18").Value = 2396
$ws.Range("D18").Value = 0.1636915537528694
$ws.Range("E18").Value = 3.691554218530655
$ws.Range("F18").Value = 2396
$ws.Range("G18").Value = 0.1124671795405447
$ws.Range("H18").Value = 0.9492811446543783
$ws.Range("I18").Value = 2.191557149635628
$ws.Range("J18").Value = 0.3017278488259763
$ws.Range("K18").Value = 0.03897371119819582
$ws.Range("E19").Value = 6.902928543509915

$wb.Save()
